# Apply the "6team-IndividualWorkSheet" update to the first worksheet (윤다은 / Table1).
# The sheet tracks a personal task list; this edit:
#   - tweaks several existing rows' result/problem notes
#   - fills in previously-blank cells in row 11 (history view row)
#   - appends two brand-new rows (16, 17) about email-validation / Korean-encoding fixes
#   - grows Table1 from A1:F15 to A1:F17 to cover the new rows

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Row 4: result note simplified to "완료" ---
$ws.Range("E4").Value = "완료"

# --- Row 7: page-merge task now finished, with a follow-up note in "문제점" ---
$ws.Range("A7").Value = "페이지 통합"
$ws.Range("B7").Value = "소비자 공급자 분리된 페이지를 하나로 합치기 (메인, 로그인, 로그아웃, 개인정보 수정, 방정보 조회)"
$ws.Range("E7").Value = "완료"
$ws.Range("F7").Value = "6.18 -> 공급자와 소비자의 reserve history를 하나의 파일로 합침"

# --- Row 10: wording tweak on content/problem columns ---
$ws.Range("B10").Value = "공급자 메인 페이지의 신청내역에서 상세 신청내역까지 구현 , 사용자의 예약 내역 (승인된것 제외) 보기 및 취소"
$ws.Range("F10").Value = "거절당한 예약신청을 직접 삭제해야함 -> 전달도 해주고 삭제도 자동으로 하려면 어떻게 해야할지 고민.."

# --- Row 11: previously-empty history-view row gets filled in ---
$ws.Range("B11").Value = "이전까지 승낙된 예약내역을 보여줌"
$ws.Range("D11").Value = 43634
$ws.Range("E11").Value = "완료"
$ws.Range("F11").Value = "날짜별로 예약 구분은 안됨.. "

# --- Row 12: wording tweak ---
$ws.Range("B12").Value = "다른 사용자의 개인정보를 열람할 수 있는 페이지를 만듬 -> 예약을 만드는 경우와 수락하는 경우에 참고, 예약된 후 서로 의사소통하기위해 전화번호를 넣었음 -> 채팅기능있으면 필요없음"
$ws.Range("E12").Value = "다른 사용자의 정보 열람 가능"
$ws.Range("F12").Value = "개인정보 유출의 문제점이 있으므로 채팅을 구현한다면 전화번호 삭제할 예정"

# --- Row 13: wording tweak ---
$ws.Range("B13").Value = "공급자에게 온 예약신청을 수락/거절 하는 기능"
$ws.Range("E13").Value = "consumer가 신청한 예약내역을 승인/거절하는 기능 구현"
$ws.Range("F13").Value = "jsp페이지를 비효율적으로 2개를 만들어서 처리하였음.. 페이지 개수 줄이는 방법이 있을듯"

# --- Row 14: result note replaced with the DB foreign-key fix ---
$ws.Range("E14").Value = "table foreign key 수정, 추가로 필요한 column 추가"

# --- Row 15: new login-session task replaces the old "페이지 통합" content that moved to row 7 ---
$ws.Range("A15").Value = "로그인 세션 문제 해결"
$ws.Range("B15").Value = "같은 세션에서 재로그인 가능 -> 기존 화면 바뀌는 문제해결 => 세션에 로그인 중일 경우 중복 로그인 불가, 그냥 메인페이지로 접속되게"
$ws.Range("E15").Value = "기존에 로그인 되어있으면 로그인 할 수 없게 수정"

# --- New row 16 (Korean-encoding fix on sending) : copy formatting from row 9 (same ht=34.5 shape) ---
$ws.Range("A9:F9").Copy()
$ws.Range("A16:F16").PasteSpecial(-4122)
$ws.Rows.Item(16).RowHeight = 34.5
$ws.Range("A16").Value = "한글 전송 깨짐 해결"
$ws.Range("B16").Value = "한글 post로 전송시 깨짐"
$ws.Range("C16").Value = 43631
$ws.Range("D16").Value = 43631
$ws.Range("E16").Value = "받는 코드에서 인코딩 바꿔서 해결!"

# --- New row 17 (signup email/password validation) : copy formatting from row 15 (same ht=51.75 shape) ---
$ws.Range("A15:F15").Copy()
$ws.Range("A17:F17").PasteSpecial(-4122)
$ws.Rows.Item(17).RowHeight = 51.75
$ws.Range("A17").Value = "회원가입시 이메일 양식 확인, 비밀번호 확인"
$ws.Range("B17").Value = "비밀번호 8자에서 15자 사이, 특수문자, 숫자 반드시 포함"
$ws.Range("C17").Value = 43630
$ws.Range("D17").Value = 43630
$ws.Range("E17").Value = "이메일 양식 확인하지않으면 다음 단계로 넘어가지못함, 비밀번호 확인도 마찬가지"
$ws.Range("F17").Value = "이메일을 확인 메일을 보내서 하면 좋을듯"

$excel.CutCopyMode = $false

# --- Grow Table1 (ListObject) to cover the two new rows so ref/autoFilter become A1:F17 ---
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:F17"))

# --- Restore the view: scrolled down a bit, with E19 selected (as in the saved file) ---
$ws.Range("E19").Select()
